# Updates the "Price" (column D) values and a few coin rows for the
# cryptos.xlsx symbol-list refresh (GitHub Actions run on
# Fri Dec 23 04:57:07 UTC 2022).
#
# All Price cells in this sheet are stored as literal text (not numbers),
# so trailing / significant zeros must be preserved exactly as strings
# (e.g. "0.03000", "0.0004000"). Setting a numeric-looking string via
# .Value normally gets auto-converted to a real number by Excel, which
# would silently strip such trailing zeros. To avoid that we force the
# cell to Text format before writing the value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param(
        [string]$Address,
        [string]$Text
    )
    $cell = $ws.Range($Address)
    $cell.NumberFormat = "@"
    $cell.Value = $Text
}

# ---- Column D (Price) updates ----
Set-TextValue "D2"  "246.34"
Set-TextValue "D3"  "22.04"
Set-TextValue "D4"  "5.448"
Set-TextValue "D5"  "0.05784"
Set-TextValue "D7"  "6.316"
Set-TextValue "D8"  "0.8182"
Set-TextValue "D9"  "0.9758"
Set-TextValue "D11" "0.07502"
Set-TextValue "D12" "0.03119"
Set-TextValue "D13" "0.03000"
Set-TextValue "D14" "4.150"
Set-TextValue "D15" "0.09409"
Set-TextValue "D16" "0.001600"
Set-TextValue "D17" "0.04815"
Set-TextValue "D19" "0.006189"
Set-TextValue "D21" "0.0009976"
Set-TextValue "D23" "3.767"
Set-TextValue "D24" "2.214"
Set-TextValue "D27" "0.0004000"
Set-TextValue "D40" "0.03889"

# ---- E9: FTXToken now flagged "Best in 24h" ----
$ws.Range("E9").Value = "8FTXTokenFTTBestin24h"

# ---- Rows 41-43 rotate (KickToken / BKEXToken / CEJI reorder) ----
# Row 41 becomes BKEXToken (was KickToken)
$ws.Range("B41").Value = "BKEXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextValue "D41" "0.1077"
$ws.Range("E41").Value = "40BKEXTokenBKK"

# Row 42 becomes CEJI (was BKEXToken)
$ws.Range("B42").Value = "CEJI"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
Set-TextValue "D42" "0.002629"
$ws.Range("E42").Value = "41CEJICEJI"

# Row 43 becomes KickToken (was CEJI)
$ws.Range("B43").Value = "KickToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
Set-TextValue "D43" "0.003069"
$ws.Range("E43").Value = "42KickTokenKICK"

# ---- Remaining column D (Price) updates ----
Set-TextValue "D44" "0.006696"
Set-TextValue "D47" "0.3801"
Set-TextValue "D49" "0.00002101"
